$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8848228454589844
$ws.Range("B1").Value = 1.671760439872742
$ws.Range("C1").Value = 4.001961708068848
$ws.Range("D1").Value = 3.781366109848022
$ws.Range("E1").Value = 0.5807180404663086
